$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.750268
$ws.Range("H2").Value = 59.250804
$ws.Range("I2").Value = 0.2138888518073023
$ws.Range("J2").Value = 0.2138888518073023
$ws.Range("M2").Value = 0.6340813333333334
$ws.Range("N2").Value = 1.902244
$ws.Range("O2").Value = 0.003487630722257058
$ws.Range("P2").Value = 0.003487630722257058
$ws.Range("Q2").Value = 12.52327626713067
$ws.Range("R2").Value = 112.709486404176
$ws.Range("S2").Value = 0.0007459653307114347
$ws.Range("T2").Value = 0.0007459653307114346
$ws.Range("G3").Value = 19.750268
$ws.Range("H3").Value = 59.250804
$ws.Range("I3").Value = 0.2138888518073023
$ws.Range("J3").Value = 0.2138888518073023
$ws.Range("O3").Value = 0.8644503444376447
$ws.Range("P3").Value = 0.8644503444376448
$ws.Range("Q3").Value = 3104.041495426122
$ws.Range("R3").Value = 27936.3734588351
$ws.Range("S3").Value = 0.1848962916161948
$ws.Range("T3").Value = 0.1848962916161949
$ws.Range("G4").Value = 19.750268
$ws.Range("H4").Value = 59.250804
$ws.Range("I4").Value = 0.2138888518073023
$ws.Range("J4").Value = 0.2138888518073023
$ws.Range("M4").Value = 24.01001466666667
$ws.Range("N4").Value = 72.030044
$ws.Range("O4").Value = 0.1320620248400982
$ws.Range("P4").Value = 0.1320620248400982
$ws.Range("Q4").Value = 474.2042243505973
$ws.Range("R4").Value = 4267.838019155376
$ws.Range("S4").Value = 0.02824659486039603
$ws.Range("T4").Value = 0.02824659486039603
$ws.Range("I5").Value = 0.6395228081370402
$ws.Range("J5").Value = 0.6395228081370402
$ws.Range("M5").Value = 0.6340813333333334
$ws.Range("N5").Value = 1.902244
$ws.Range("O5").Value = 0.003487630722257058
$ws.Range("P5").Value = 0.003487630722257058
$ws.Range("Q5").Value = 37.44431155601689
$ws.Range("R5").Value = 336.998804004152
$ws.Range("S5").Value = 0.002230419393242848
$ws.Range("T5").Value = 0.002230419393242848
$ws.Range("I6").Value = 0.6395228081370402
$ws.Range("J6").Value = 0.6395228081370402
$ws.Range("O6").Value = 0.8644503444376447
$ws.Range("P6").Value = 0.8644503444376448
$ws.Range("R6").Value = 83529.12203042342
$ws.Range("S6").Value = 0.5528357117697942
$ws.Range("T6").Value = 0.5528357117697942
$ws.Range("I7").Value = 0.6395228081370402
$ws.Range("J7").Value = 0.6395228081370402
$ws.Range("M7").Value = 24.01001466666667
$ws.Range("N7").Value = 72.030044
$ws.Range("O7").Value = 0.1320620248400982
$ws.Range("P7").Value = 0.1320620248400982
$ws.Range("Q7").Value = 1417.859858635172
$ws.Range("R7").Value = 12760.73872771655
$ws.Range("S7").Value = 0.08445667697400314
$ws.Range("T7").Value = 0.08445667697400314
$ws.Range("G8").Value = 13.53581066666667
$ws.Range("H8").Value = 40.607432
$ws.Range("I8").Value = 0.1465883400556574
$ws.Range("J8").Value = 0.1465883400556574
$ws.Range("M8").Value = 0.6340813333333334
$ws.Range("N8").Value = 1.902244
$ws.Range("O8").Value = 0.003487630722257058
$ws.Range("P8").Value = 0.003487630722257058
$ws.Range("Q8").Value = 8.582804875267557
$ws.Range("R8").Value = 77.245243877408
$ws.Range("S8").Value = 0.0005112459983027757
$ws.Range("T8").Value = 0.0005112459983027757
$ws.Range("G9").Value = 13.53581066666667
$ws.Range("H9").Value = 40.607432
$ws.Range("I9").Value = 0.1465883400556574
$ws.Range("J9").Value = 0.1465883400556574
$ws.Range("O9").Value = 0.8644503444376447
$ws.Range("P9").Value = 0.8644503444376448
$ws.Range("Q9").Value = 2127.349258428537
$ws.Range("R9").Value = 19146.14332585684
$ws.Range("S9").Value = 0.1267183410516557
$ws.Range("T9").Value = 0.1267183410516557
$ws.Range("G10").Value = 13.53581066666667
$ws.Range("H10").Value = 40.607432
$ws.Range("I10").Value = 0.1465883400556574
$ws.Range("J10").Value = 0.1465883400556574
$ws.Range("M10").Value = 24.01001466666667
$ws.Range("N10").Value = 72.030044
$ws.Range("O10").Value = 0.1320620248400982
$ws.Range("P10").Value = 0.1320620248400982
$ws.Range("Q10").Value = 324.9950126318898
$ws.Range("R10").Value = 2924.955113687008
$ws.Range("S10").Value = 0.01935875300569899
$ws.Range("T10").Value = 0.01935875300569899
